# Add MASE measurement to GA
# - Rename the "VAR" sheet to "Salmon" and make it the active/selected sheet
# - Update the GA rank-selection data on that sheet: drop the ARIMA column
#   entry (G1/G2) and swap in the new MASE-driven rank values in F3/F4
# - The dependent COUNTIF() formulas in column C recalc automatically

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("VAR")
$ws.Name = "Salmon"

# Clear the "ARIMA" header/value now that ARIMA is no longer tracked here
$ws.Range("G1").ClearContents()
$ws.Range("G2").ClearContents()

# New GA-selected rank values (MASE measurement added to the run)
$ws.Range("F3").Value = 59
$ws.Range("F4").Value = 1

# Leave the sheet with G6 selected and activated (becomes the active tab)
$ws.Range("G6").Select() | Out-Null
